$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("N2").Value = 70831.95557958097
$ws.Range("O2").Value = 69610.44223910036

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 47338.61932520662
$ws.Range("I2").Value = 283167.7766510943
$ws.Range("L2").Value = 178095.3756971828
$ws.Range("M2").Value = 114008.3253427963
$ws.Range("N2").Value = 33931.82461160053
$ws.Range("O2").Value = 50657.26889981552

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 29273.60317916481
$ws.Range("B2").Value = 22330.72247668595
$ws.Range("E2").Value = 110739.3594843864
$ws.Range("I2").Value = 150386.9441391908
$ws.Range("M2").Value = 35556.98862372932
$ws.Range("N2").Value = 44813.41193308897
$ws.Range("O2").Value = 26775.55841092002

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 1041.156112142704

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 34256.25291363284
$ws.Range("N2").Value = 5271.89502409355
$ws.Range("O2").Value = 22972.54525065989
